# Test data added to excel file.
# Builds a 3-column "Test case" table (TC_ID / UserName / Password header rows
# plus one sample row) on Sheet1, matching the authored TestData.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell values -----------------------------------------------------------
# Written in this order so the shared-strings table comes out in the same
# order as the authored workbook (TC_01, UserName, Password, TC_ID, student,
# Password123, Verify...).
$ws.Range("A2").Value = "TC_01"
$ws.Range("B2").Value = "UserName"
$ws.Range("C2").Value = "Password"
$ws.Range("A1").Value = "TC_ID"
$ws.Range("B3").Value = "student"
$ws.Range("C3").Value = "Password123"
$ws.Range("A3").Value = "Verify that user is able to login with valid credentials."

# --- header row (A1:C1) -----------------------------------------------------
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").Interior.Color = 255
$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A1:C1").Merge()

# --- column-title row (A2:C2) ----------------------------------------------
$ws.Range("A2:C2").Font.Bold = $true
$ws.Range("A2:C2").Interior.ThemeColor = 10

# --- data row (B3:C3) -------------------------------------------------------
$ws.Range("B3:C3").Font.Name = "Courier New"
$ws.Range("B3:C3").Font.Family = 3
$ws.Range("B3:C3").Font.Size = 14
$ws.Range("B3:C3").Font.Color = 16711722
$ws.Range("B3:C3").VerticalAlignment = -4108

$ws.Rows.Item(3).RowHeight = 18

# --- column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 46
$ws.Columns.Item(2).ColumnWidth = 30
$ws.Columns.Item(3).ColumnWidth = 30.8

# --- sheet view / page setup -------------------------------------------------
$ws.Range("B8").Select()
$ws.PageSetup.Orientation = 1
